# Update project list - clean up trailing commas in "Assigned members" and
# refresh/normalize "Due date" values (folders were updated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E: Assigned members - strip trailing ", " ---
$ws.Range("E2").Value = "admin, HoanLeader"
$ws.Range("E7").Value = "Nobita"
$ws.Range("E8").Value = "admin, HoanLeader"
$ws.Range("E10").Value = "admin, HoanLeader, HoanTester"
$ws.Range("E12").Value = "HoanTester"
$ws.Range("E14").Value = "HoanLeader"

# --- Column D: Due date updates ---
$ws.Range("D3").Value = "Friday, 30 August 2019"
$ws.Range("D4").Value = "Sunday, 30 June 2019"
$ws.Range("D5").Value = "Sunday, 30 June 2019"
$ws.Range("D6").Value = "Sunday, 30 June 2019"
$ws.Range("D7").Value = "Sunday, 30 June 2019"
$ws.Range("D8").Value = "Tuesday, 30 July 2019"
$ws.Range("D9").Value = "Sunday, 30 June 2019"
$ws.Range("D10").Value = "Sunday, 30 June 2019"
$ws.Range("D11").Value = "Sunday, 30 June 2019"
$ws.Range("D12").Value = "Sunday, 30 June 2019"
$ws.Range("D13").Value = "Sunday, 30 June 2019"
$ws.Range("D14").Value = "Sunday, 30 June 2019"
$ws.Range("D15").Value = "Monday, 17 June 2019"

# --- Column widths (bestFit was recalculated by Excel after the content change) ---
$ws.Columns.Item(4).ColumnWidth = 24.666666666666668
$ws.Columns.Item(5).ColumnWidth = 34.333333333333336
